# Add: screenshots of project and updated rubric
# - Lower the "Execution" score from 3 to 2 (F4); the Total Score formula
#   in F12 (=SUM(F4:F8)) recalculates automatically from 15 to 14.
# - Move the active selection to F4 (instead of C11) and widen the workbook
#   window to better show the updated rubric/screenshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F4").Value = 2

$ws.Range("F4").Select()
$excel.ActiveWindow.Width = 22665
